# Mise à jour du classement - 26.03.2025 à 18:00
$wb = $excel.ActiveWorkbook

# --- Sheet "leaderboard2" : "Qui a attrapé le plus de Cobblemons ?" ---
$ws2 = $wb.Worksheets.Item("leaderboard2")
$ws2.Range("D3").Value = 101
$ws2.Range("D4").Value = 51
$ws2.Range("D5").Value = 21
$ws2.Range("D6").Value = 8
$ws2.Range("B13").Value = "Dernière update le 26.03.25 à 18:00"

# --- Sheet "leaderboard3" : "Qui a attrapé le plus de Shiny Cobblemons ?" ---
# Leaderboard re-ordered: BKZRackham now 1st, ArtyumsM 2nd, Lokys 3rd.
$ws3 = $wb.Worksheets.Item("leaderboard3")
$ws3.Range("C3").Value = "BKZRackham"
$ws3.Range("D3").Value = 13
$ws3.Range("C4").Value = "ArtyumsM"
$ws3.Range("D4").Value = 11
$ws3.Range("C5").Value = "Lokys"
$ws3.Range("D5").Value = 8
$ws3.Range("B13").Value = "Dernière update le 26.03.25 à 18:00"
